$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update parameter lower/upper bounds per the model-equation revision
$ws.Range("C5").Value = 0.0001
$ws.Range("B6").Value = 0.001
$ws.Range("B7").Value = 0.1
$ws.Range("C8").Value = 10
$ws.Range("B12").Value = 0.0001
$ws.Range("B19").Value = 0.01
$ws.Range("C19").Value = 10
$ws.Range("B22").Value = 0.0001
$ws.Range("C22").Value = 1

# Move the active selection to C22 (and let scrolling reset to the top)
$ws.Activate()
$ws.Range("C22").Select()
